$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# ColumnWidth (characters) maps to the stored XML "width" as width = ColumnWidth + 5/6,
# quantized to the nearest 1/6th of a character by the engine (same pixel-grid snapping
# real Excel performs). Subtract 5/6 from the desired stored width so the result lands
# as close as possible to the target width from the diff.
$ws.Columns.Item(6).ColumnWidth = 20.25 - (5/6)    # F: 13.5 -> 20.25
$ws.Columns.Item(8).ColumnWidth = 10.8 - (5/6)     # H: 32.4 -> 10.8
$ws.Columns.Item(9).ColumnWidth = 37.8 - (5/6)     # I: 28.35 -> 37.8

# --- Row 3 updates ---
$ws.Range("B3").Value = "16-12-2025"
$ws.Range("C3").Value = "Nithin"
$ws.Range("D3").Value = 8943053925
$ws.Range("E3").Value = "'03-01-2026"
$ws.Range("F3").Value = "RAYAN K B"
$ws.Range("H3").Value = "ENQUIRY"
$ws.Range("I3").Value = "ENQUIRY WITHOUT BRIDE/FAMILY"
$ws.Range("K3").Value = "will update with family"

# --- Row 4 updates ---
$ws.Range("B4").Value = "16-12-2025"
$ws.Range("C4").Value = "adharv"
$ws.Range("D4").Value = 7994499796
$ws.Range("E4").Value = "23-12-2025"
$ws.Range("H4").Value = "ENQUIRY"
$ws.Range("I4").Value = "Enquiry for Relative/Friend"
$ws.Range("K4").Value = "will update tomorrow"

# --- Row 5 updates ---
$ws.Range("B5").Value = "16-12-2025"
$ws.Range("C5").Value = "joyal"
$ws.Range("D5").Value = 9567564713
$ws.Range("E5").Value = "'08-01-2026"
$ws.Range("F5").Value = "ATHULKIRSHNA CS"
$ws.Range("H5").Value = "PRODUCT"
$ws.Range("I5").Value = "PRODUCT NOT AVAILABLE"
$ws.Range("K5").Value = "navi blue suit not available in store"

# --- Row 6 updates ---
$ws.Range("B6").Value = "16-12-2025"
$ws.Range("C6").Value = "shibin"
$ws.Range("D6").Value = 7907275586
$ws.Range("E6").Value = "16-12-2025"
$ws.Range("F6").Value = "ATHULKIRSHNA CS"
$ws.Range("K6").Value = "beach sand colour suit not available in store"

# --- Row 7 (new row) ---
$ws.Range("A7").Value = 5
$ws.Range("A7").NumberFormat = "0"
$ws.Range("B7").Value = "16-12-2025"
$ws.Range("C7").Value = "SAHAL"
$ws.Range("D7").Value = 7907034399
$ws.Range("D7").NumberFormat = "0"
$ws.Range("E7").Value = "25-01-2026"
$ws.Range("F7").Value = "RASEEB E A"
$ws.Range("G7").Value = "Loss"
$ws.Range("H7").Value = "PRODUCT"
$ws.Range("I7").Value = "PRODUCT NOT AVAILABLE"
$ws.Range("J7").Value = "-"
$ws.Range("K7").Value = "PRODUCT WAS NOT AVAILABLE FOR TRAIL"
